# Insert two new data rows right before the existing row 40.
# This pushes the current rows 40..160 down to 42..162, matching the
# target workbook's new dimension of A1:R162.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40:A41").EntireRow.Insert()

# New row 40
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44487
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = 100112039
$ws.Range("G40").Value = "Ciboulette"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 105
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = 5619
$ws.Range("N40").Value = "$/docena de atados"
$ws.Range("O40").Value = "Provincia de Cautín"
$ws.Range("P40").Value = 1873
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"

# New row 41
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 44487
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112039
$ws.Range("G41").Value = "Ciboulette"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 75
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 2000
$ws.Range("M41").Value = 2000
$ws.Range("N41").Value = "$/docena de atados"
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 667
$ws.Range("Q41").Value = 3
$ws.Range("R41").Value = "Hortaliza"
